$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.479.88"
$ws.Range("D3").Value = "2.325.18"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "511.59"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.47"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("E10").Value = "  -0.20%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "2.741.65"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.53"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "56.482.41"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "2.326.10"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.41"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "327.44"
$ws.Range("E19").Value = "  +1.23%  "
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.71"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.60"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.60"
$ws.Range("E24").Value = "  +9.15%  "
$ws.Range("E25").Value = "  -0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  +2.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "167.81"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.67"
$ws.Range("E29").Value = "  -3.68%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0718"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.35"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.25"
$ws.Range("E35").Value = "  +0.64%  "
$ws.Range("E36").Value = "  -2.70%  "
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.49"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "149.04"
$ws.Range("E40").Value = "  +8.43%  "
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "277.76"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.02"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.14"
$ws.Range("E48").Value = "  +1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.379"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("E51").Value = "  +1.81%  "
